# Add the new "S*-unmerged" / "S*-unmergedND" variant columns/rows that a new
# ChangeKeysDict-based run contributed to the results.

$wb = $excel.ActiveWorkbook

# --- Sheet "results": insert a "S*-unmerged" column (after S*-MM0) and an
#     "S*-unmergedND" column (at the end), then fill in their values. ---
$ws1 = $wb.Worksheets.Item("results")
$ws1.Columns("G").Insert()
$ws1.Columns("M").Insert()
$ws1.Range("G1").Value = "S*-unmerged"
$ws1.Range("G2").Value = 114
$ws1.Range("M1").Value = "S*-unmergedND"
$ws1.Range("M2").Value = $false
$ws1.Range("J2").Value = $false

# --- Sheet "stats": insert a "S*-unmerged" row into each of the two
#     run-0 / Average blocks (just before their "Kruskal" summary row). ---
$ws2 = $wb.Worksheets.Item("stats")

# Insert the two new rows (extends merges A2:A6->A2:A7 and A7:A11->A8:A13 automatically)
$ws2.Rows("6").Insert()
$ws2.Rows("12").Insert()

# Fix formatting on the two brand-new rows (copy style from the row above) for columns A and B
$ws2.Range("A5:B5").Copy() | Out-Null
$ws2.Range("A6:B6").PasteSpecial(-4122) | Out-Null
$ws2.Range("A11:B11").Copy() | Out-Null
$ws2.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Now set every data cell in rows 2-13 explicitly to the final values
$ws2.Range("A2").Value = "run 0"
$ws2.Range("B2").Value = "S*-BS"
$ws2.Range("C2").Value = 46
$ws2.Range("D2").Value = 0.00008912896737456322
$ws2.Range("E2").Value = 0.0222397712059319
$ws2.Range("F2").Value = 46
$ws2.Range("G2").Value = 0.001979269552975893
$ws2.Range("H2").Value = 0.003122320864349604
$ws2.Range("I2").Value = 0.005609294399619102
$ws2.Range("J2").Value = 0.009039589669555426
$ws2.Range("K2").Value = 0.0007170438766479492
$ws2.Range("A3").ClearContents() | Out-Null
$ws2.Range("B3").Value = "S*-HS"
$ws2.Range("C3").Value = 46
$ws2.Range("D3").Value = 0.001718624960631132
$ws2.Range("E3").Value = 0.02583552291616797
$ws2.Range("F3").Value = 46
$ws2.Range("G3").Value = 0.001975357532501221
$ws2.Range("H3").Value = 0.004766576923429966
$ws2.Range("I3").Value = 0.005655502434819937
$ws2.Range("J3").Value = 0.01113275717943907
$ws2.Range("K3").Value = 0.0006633377633988857
$ws2.Range("A4").ClearContents() | Out-Null
$ws2.Range("B4").Value = "S*-MM"
$ws2.Range("C4").Value = 46
$ws2.Range("D4").Value = 0.00203008996322751
$ws2.Range("E4").Value = 0.02835799194872379
$ws2.Range("F4").Value = 46
$ws2.Range("G4").Value = 0.002220173366367817
$ws2.Range("H4").Value = 0.005485605914145708
$ws2.Range("I4").Value = 0.006355251651257277
$ws2.Range("J4").Value = 0.01170287234708667
$ws2.Range("K4").Value = 0.0007458264008164406
$ws2.Range("A5").ClearContents() | Out-Null
$ws2.Range("B5").Value = "S*-MM0"
$ws2.Range("C5").Value = 46
$ws2.Range("D5").Value = 0.000132476445287466
$ws2.Range("E5").Value = 0.02243422856554389
$ws2.Range("F5").Value = 46
$ws2.Range("G5").Value = 0.001913452055305243
$ws2.Range("H5").Value = 0.003099545370787382
$ws2.Range("I5").Value = 0.005960895679891109
$ws2.Range("J5").Value = 0.008982968982309103
$ws2.Range("K5").Value = 0.0007380195893347263
$ws2.Range("A6").ClearContents() | Out-Null
$ws2.Range("B6").Value = "S*-unmerged"
$ws2.Range("C6").Value = 85
$ws2.Range("D6").Value = 0.002546215895563364
$ws2.Range("E6").Value = 0.06322001619264483
$ws2.Range("F6").Value = 85
$ws2.Range("G6").Value = 0.003238909412175417
$ws2.Range("H6").Value = 0.008396383840590715
$ws2.Range("I6").Value = 0.03092650882899761
$ws2.Range("J6").Value = 0.01608820864930749
$ws2.Range("K6").Value = 0.001378839369863272
$ws2.Range("A7").ClearContents() | Out-Null
$ws2.Range("B7").Value = "Kruskal"
$ws2.Range("C7").Value = 1425
$ws2.Range("D7").ClearContents() | Out-Null
$ws2.Range("E7").Value = 0.01139114506077021
$ws2.Range("F7").ClearContents() | Out-Null
$ws2.Range("G7").ClearContents() | Out-Null
$ws2.Range("H7").ClearContents() | Out-Null
$ws2.Range("I7").ClearContents() | Out-Null
$ws2.Range("J7").ClearContents() | Out-Null
$ws2.Range("K7").ClearContents() | Out-Null
$ws2.Range("A8").Value = "Average"
$ws2.Range("B8").Value = "S*-BS"
$ws2.Range("C8").Value = 46
$ws2.Range("D8").Value = 0.00008912896737456322
$ws2.Range("E8").Value = 0.0222397712059319
$ws2.Range("F8").Value = 46
$ws2.Range("G8").Value = 0.001979269552975893
$ws2.Range("H8").Value = 0.003122320864349604
$ws2.Range("I8").Value = 0.005609294399619102
$ws2.Range("J8").Value = 0.009039589669555426
$ws2.Range("K8").Value = 0.0007170438766479492
$ws2.Range("A9").ClearContents() | Out-Null
$ws2.Range("B9").Value = "S*-HS"
$ws2.Range("C9").Value = 46
$ws2.Range("D9").Value = 0.001718624960631132
$ws2.Range("E9").Value = 0.02583552291616797
$ws2.Range("F9").Value = 46
$ws2.Range("G9").Value = 0.001975357532501221
$ws2.Range("H9").Value = 0.004766576923429966
$ws2.Range("I9").Value = 0.005655502434819937
$ws2.Range("J9").Value = 0.01113275717943907
$ws2.Range("K9").Value = 0.0006633377633988857
$ws2.Range("A10").ClearContents() | Out-Null
$ws2.Range("B10").Value = "S*-MM"
$ws2.Range("C10").Value = 46
$ws2.Range("D10").Value = 0.00203008996322751
$ws2.Range("E10").Value = 0.02835799194872379
$ws2.Range("F10").Value = 46
$ws2.Range("G10").Value = 0.002220173366367817
$ws2.Range("H10").Value = 0.005485605914145708
$ws2.Range("I10").Value = 0.006355251651257277
$ws2.Range("J10").Value = 0.01170287234708667
$ws2.Range("K10").Value = 0.0007458264008164406
$ws2.Range("A11").ClearContents() | Out-Null
$ws2.Range("B11").Value = "S*-MM0"
$ws2.Range("C11").Value = 46
$ws2.Range("D11").Value = 0.000132476445287466
$ws2.Range("E11").Value = 0.02243422856554389
$ws2.Range("F11").Value = 46
$ws2.Range("G11").Value = 0.001913452055305243
$ws2.Range("H11").Value = 0.003099545370787382
$ws2.Range("I11").Value = 0.005960895679891109
$ws2.Range("J11").Value = 0.008982968982309103
$ws2.Range("K11").Value = 0.0007380195893347263
$ws2.Range("A12").ClearContents() | Out-Null
$ws2.Range("B12").Value = "S*-unmerged"
$ws2.Range("C12").Value = 85
$ws2.Range("D12").Value = 0.002546215895563364
$ws2.Range("E12").Value = 0.06322001619264483
$ws2.Range("F12").Value = 85
$ws2.Range("G12").Value = 0.003238909412175417
$ws2.Range("H12").Value = 0.008396383840590715
$ws2.Range("I12").Value = 0.03092650882899761
$ws2.Range("J12").Value = 0.01608820864930749
$ws2.Range("K12").Value = 0.001378839369863272
$ws2.Range("A13").ClearContents() | Out-Null
$ws2.Range("B13").Value = "Kruskal"
$ws2.Range("C13").Value = 1425
$ws2.Range("D13").ClearContents() | Out-Null
$ws2.Range("E13").Value = 0.01139114506077021
$ws2.Range("F13").ClearContents() | Out-Null
$ws2.Range("G13").ClearContents() | Out-Null
$ws2.Range("H13").ClearContents() | Out-Null
$ws2.Range("I13").ClearContents() | Out-Null
$ws2.Range("J13").ClearContents() | Out-Null
$ws2.Range("K13").ClearContents() | Out-Null
